# Auto-generated edit script
# Applies updated '最低票价' (min ticket price) values to columns G (as text,
# displaying price/100 or '不可售' when not sellable), and a handful of
# incidental '想去人数' (want-to-go count) bumps in column F, across the four
# worksheets, matching the upstream scrape refresh described in the diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")

$ws1Changes = @(
    @{ Row = 2; F = $null; G = "88" },
    @{ Row = 3; F = $null; G = "80" },
    @{ Row = 4; F = $null; G = "458" },
    @{ Row = 5; F = $null; G = "不可售" },
    @{ Row = 6; F = 509; G = "不可售" },
    @{ Row = 7; F = $null; G = "598" },
    @{ Row = 8; F = 9029; G = "8" },
    @{ Row = 9; F = $null; G = "598" },
    @{ Row = 10; F = $null; G = "68" },
    @{ Row = 11; F = $null; G = "75" },
    @{ Row = 12; F = $null; G = "158" },
    @{ Row = 13; F = $null; G = "49" },
    @{ Row = 14; F = 2283; G = "75" },
    @{ Row = 15; F = $null; G = "70" },
    @{ Row = 16; F = 3564; G = "19" },
    @{ Row = 17; F = 242; G = "60" },
    @{ Row = 18; F = $null; G = "21" },
    @{ Row = 19; F = 99; G = "48" },
    @{ Row = 20; F = $null; G = "55" },
    @{ Row = 21; F = $null; G = "65" },
    @{ Row = 22; F = $null; G = "75" },
    @{ Row = 23; F = $null; G = "60" },
    @{ Row = 24; F = $null; G = "93" },
    @{ Row = 25; F = $null; G = "60" },
    @{ Row = 26; F = $null; G = "29" },
    @{ Row = 27; F = $null; G = "70" },
    @{ Row = 28; F = $null; G = "19" },
    @{ Row = 29; F = $null; G = "60" },
    @{ Row = 30; F = $null; G = "75" },
    @{ Row = 31; F = $null; G = "不可售" },
    @{ Row = 32; F = $null; G = "75" }
)

foreach ($change in $ws1Changes) {
    if ($null -ne $change.F) {
        $ws1.Cells.Item($change.Row, 6).Value = $change.F
    }
    $gCell = $ws1.Cells.Item($change.Row, 7)
    if ($change.G -match "^[0-9]+$") {
        $gCell.NumberFormat = "@"
    }
    $gCell.Value = $change.G
}

$ws2 = $wb.Worksheets.Item("演出")

$ws2Changes = @(
    @{ Row = 2; F = $null; G = "380" },
    @{ Row = 3; F = $null; G = "90" },
    @{ Row = 4; F = $null; G = "180" },
    @{ Row = 5; F = $null; G = "88" },
    @{ Row = 6; F = $null; G = "680" }
)

foreach ($change in $ws2Changes) {
    if ($null -ne $change.F) {
        $ws2.Cells.Item($change.Row, 6).Value = $change.F
    }
    $gCell = $ws2.Cells.Item($change.Row, 7)
    if ($change.G -match "^[0-9]+$") {
        $gCell.NumberFormat = "@"
    }
    $gCell.Value = $change.G
}

$ws3 = $wb.Worksheets.Item("本地生活")

$ws3Changes = @(
    @{ Row = 2; F = 168; G = "30" },
    @{ Row = 3; F = $null; G = "20" }
)

foreach ($change in $ws3Changes) {
    if ($null -ne $change.F) {
        $ws3.Cells.Item($change.Row, 6).Value = $change.F
    }
    $gCell = $ws3.Cells.Item($change.Row, 7)
    if ($change.G -match "^[0-9]+$") {
        $gCell.NumberFormat = "@"
    }
    $gCell.Value = $change.G
}

$ws4 = $wb.Worksheets.Item("全部类型")

$ws4Changes = @(
    @{ Row = 2; F = 168; G = "30" },
    @{ Row = 3; F = $null; G = "88" },
    @{ Row = 4; F = $null; G = "20" },
    @{ Row = 5; F = $null; G = "80" },
    @{ Row = 6; F = $null; G = "380" },
    @{ Row = 7; F = $null; G = "90" },
    @{ Row = 8; F = $null; G = "458" },
    @{ Row = 9; F = $null; G = "不可售" },
    @{ Row = 10; F = 509; G = "不可售" },
    @{ Row = 11; F = $null; G = "598" },
    @{ Row = 12; F = 9029; G = "8" },
    @{ Row = 13; F = $null; G = "598" },
    @{ Row = 14; F = $null; G = "68" },
    @{ Row = 15; F = $null; G = "75" },
    @{ Row = 16; F = $null; G = "158" },
    @{ Row = 17; F = $null; G = "49" },
    @{ Row = 18; F = $null; G = "180" },
    @{ Row = 19; F = 2283; G = "75" },
    @{ Row = 20; F = $null; G = "70" },
    @{ Row = 21; F = 3564; G = "19" },
    @{ Row = 22; F = 242; G = "60" },
    @{ Row = 23; F = 86; G = "21" },
    @{ Row = 24; F = 99; G = "48" },
    @{ Row = 25; F = $null; G = "55" },
    @{ Row = 26; F = $null; G = "65" },
    @{ Row = 27; F = $null; G = "75" },
    @{ Row = 28; F = $null; G = "88" },
    @{ Row = 29; F = $null; G = "60" },
    @{ Row = 30; F = $null; G = "93" },
    @{ Row = 31; F = $null; G = "60" },
    @{ Row = 32; F = $null; G = "29" },
    @{ Row = 33; F = $null; G = "70" },
    @{ Row = 34; F = $null; G = "19" },
    @{ Row = 35; F = $null; G = "60" },
    @{ Row = 36; F = $null; G = "75" },
    @{ Row = 37; F = $null; G = "不可售" },
    @{ Row = 38; F = $null; G = "75" },
    @{ Row = 39; F = $null; G = "680" }
)

foreach ($change in $ws4Changes) {
    if ($null -ne $change.F) {
        $ws4.Cells.Item($change.Row, 6).Value = $change.F
    }
    $gCell = $ws4.Cells.Item($change.Row, 7)
    if ($change.G -match "^[0-9]+$") {
        $gCell.NumberFormat = "@"
    }
    $gCell.Value = $change.G
}
